$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 27.82738278199502

$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 0.04071648406533734
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 1.345027435901867

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.1494219747398047
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 5.586269137925634

$ws.Range("B5").Value = 0.0001021024915524027
$ws.Range("C5").Value = 0.002571899574220771
$ws.Range("D5").Value = 22.3905356188092
$ws.Range("E5").Value = 2195978.878461985
$ws.Range("G5").Value = 2196001.271671606

$ws.Range("B6").Value = 3.286832544864788
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 0.1494219747398047
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 5.586269137925634
